$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to show a 3-census-year comparison (1989 / 2002 / 2014) with a
# subtitle line. The new layout only keeps the most recent (2014) figure and
# drops the subtitle, so:

# 1. Remove the blank spacer row (old row 3) - rows 4,5,6 shift up to 3,4,5
$ws.Rows("3").Delete()

# 2. Remove the now-unneeded 2002/2014 columns (old columns C:D), keeping only
#    the single remaining data column (old column B)
$ws.Columns("C:D").Delete()

# 3. The single remaining data column used to hold the 1989 figures - update it
#    to hold the 2014 figures instead (value itself unchanged at 3044.5)
$ws.Range("B4").Value = 2014

# 4. Clear the subtitle row completely (text + formatting)
$ws.Range("A2:B2").Clear()

# 5. Clear the left-over blank formatted cell next to the title
$ws.Range("B1").Clear()

# 6. Rename the sheet from the generic "1" to the municipality name
$ws.Name = "მესტია"

# 7. Restore the active selection to A2, matching the saved workbook state
$ws.Range("A2").Select() | Out-Null
